$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.921.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.787.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "

# Row 7
$ws.Range("E7").Value = "  -1.92%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.79%  "

# Row 10
$ws.Range("E10").Value = "  -3.30%  "

# Row 11
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0848"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.227.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.785.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.947"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.862.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.21%  "

# Row 20
$ws.Range("E20").Value = "  -2.37%  "

# Row 21
$ws.Range("E21").Value = "  -3.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.48%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "

# Row 25
$ws.Range("E25").Value = "  -4.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.29%  "

# Row 29
$ws.Range("E29").Value = "  -1.07%  "

# Row 30
$ws.Range("B30").Value = "VeChain"
$ws.Range("C30").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0472"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.50%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.85%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0845"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.25%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("E38").Value = "  +0.66%  "

# Row 39
$ws.Range("E39").Value = "  -3.49%  "

# Row 40
$ws.Range("E40").Value = "  -4.83%  "

# Row 41
$ws.Range("E41").Value = "  +1.58%  "

# Row 42
$ws.Range("E42").Value = "  -2.23%  "

# Row 43
$ws.Range("E43").Value = "  -1.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.61%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.080.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "

# Row 47
$ws.Range("E47").Value = "  -4.80%  "

# Row 49
$ws.Range("E49").Value = "  -5.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.88%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.86%  "
